# Re-applies a row-content permutation within three blocks of the
# "Artfynd" sheet. The row *numbers* stay put; the data that lives in
# each row moves to a different row per the mapping below (captured
# from the canonical diff). We snapshot every source row's full
# A:AY content BEFORE writing anything, since several of the mappings
# form rotation cycles (e.g. 32<-33<-34<-35<-36<-37<-32) where a naive
# sequential copy would clobber a value before it's been read.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (i.e. target row's new content == source row's old content)
$mapping = @{
    30 = 31;  31 = 30;
    32 = 33;  33 = 34;  34 = 35;  35 = 36;  36 = 37;  37 = 32;
    45 = 47;  46 = 48;  47 = 45;  48 = 46;
    101 = 103; 102 = 104; 103 = 101; 104 = 102;
}

# Snapshot the full row (columns A:AY) for every row that is used as a source.
$snapshots = @{}
foreach ($src in $mapping.Values) {
    if (-not $snapshots.ContainsKey($src)) {
        $snapshots[$src] = $ws.Range("A$($src):AY$($src)").Value2
    }
}

# Now write each target row from its captured snapshot.
# Columns Y and AA hold date-like text ("2023-08-13") that must stay
# literal text (as it was in the source file) rather than be
# auto-coerced into a real date value by Excel's type inference, so we
# force those two columns to Text format on the destination row first.
foreach ($tgt in $mapping.Keys) {
    $src = $mapping[$tgt]
    $ws.Range("Y$($tgt)").NumberFormat = "@"
    $ws.Range("AA$($tgt)").NumberFormat = "@"
    $ws.Range("A$($tgt):AY$($tgt)").Value2 = $snapshots[$src]
}
